# Fix timer calculator inputs on Sheet1: desired Tloop (ms) and fclk (MHz),
# so the blink timer computes at the correct rate.
$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item("Sheet1")

# desired Tloop (ms): 22 -> 1000
$sheet1.Range("C9").Value = 1000

# fclk (MHz): was a hard-coded 50, now derived as 50/32
$sheet1.Range("D9").Formula = "=50/32"

# Sheet1 becomes the active sheet/tab, with C10 selected
$sheet1.Activate()
$sheet1.Range("C10").Select()
